$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2: sesion for Juanjo changes from 3 to 4
$ws.Range("B2").Value = 4

# Update row 4 (Daniel): sesion changes from 4 to 6, escena changes from Montaña.mp4 to Costa.mp4
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = "Costa.mp4"

# Add new row 6: DS3 | 2 | Costa.mp4
$ws.Range("A6").Value = "DS3"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "Costa.mp4"
